$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "G2" = 2.401444666666666
    "H2" = 7.204333999999999
    "I2" = 0.5723125574599716
    "J2" = 0.5723125574599716
    "K2" = 3
    "M2" = 13.441269
    "N2" = 40.323807
    "O2" = 0.08973082133481231
    "P2" = 0.08973082133481232
    "Q2" = 32.278463753282
    "R2" = 290.506173779538
    "S2" = 0.05135407584111021
    "T2" = 0.05135407584111022
    "E3" = 3
    "G3" = 2.401444666666666
    "H3" = 7.204333999999999
    "I3" = 0.5723125574599716
    "J3" = 0.5723125574599716
    "K3" = 3
    "M3" = 54.711535
    "N3" = 164.134605
    "O3" = 0.3652416280068742
    "P3" = 0.3652416280068742
    "Q3" = 131.3867239308966
    "R3" = 1182.48051537807
    "S3" = 0.2090323702154577
    "T3" = 0.2090323702154578
    "E4" = 3
    "G4" = 2.401444666666666
    "H4" = 7.204333999999999
    "I4" = 0.5723125574599716
    "J4" = 0.5723125574599716
    "K4" = 3
    "M4" = 63.67711
    "N4" = 191.03133
    "O4" = 0.4250937452800914
    "P4" = 0.4250937452800915
    "Q4" = 152.9170561982467
    "R4" = 1376.25350578422
    "S4" = 0.2432864885214868
    "T4" = 0.2432864885214869
    "E5" = 3
    "G5" = 2.401444666666666
    "H5" = 7.204333999999999
    "I5" = 0.5723125574599716
    "J5" = 0.5723125574599716
    "K5" = 3
    "M5" = 17.96553866666667
    "N5" = 53.896616
    "O5" = 0.119933805378222
    "P5" = 0.119933805378222
    "Q5" = 43.14324701486044
    "R5" = 388.289223133744
    "S5" = 0.06863962288191672
    "T5" = 0.06863962288191673
    "E6" = 3
    "G6" = 0.9802360000000001
    "H6" = 2.940708
    "I6" = 0.2336099514851752
    "J6" = 0.2336099514851752
    "K6" = 3
    "M6" = 13.441269
    "N6" = 40.323807
    "O6" = 0.08973082133481231
    "P6" = 0.08973082133481232
    "Q6" = 13.175615759484
    "R6" = 118.580541835356
    "S6" = 0.02096201281875043
    "T6" = 0.02096201281875043
    "E7" = 3
    "G7" = 0.9802360000000001
    "H7" = 2.940708
    "I7" = 0.2336099514851752
    "J7" = 0.2336099514851752
    "K7" = 3
    "M7" = 54.711535
    "N7" = 164.134605
    "O7" = 0.3652416280068742
    "P7" = 0.3652416280068742
    "Q7" = 53.63021622226
    "R7" = 482.67194600034
    "S7" = 0.08532407899905228
    "T7" = 0.0853240789990523
    "E8" = 3
    "G8" = 0.9802360000000001
    "H8" = 2.940708
    "I8" = 0.2336099514851752
    "J8" = 0.2336099514851752
    "K8" = 3
    "M8" = 63.67711
    "N8" = 191.03133
    "O8" = 0.4250937452800914
    "P8" = 0.4250937452800915
    "Q8" = 62.41859559796001
    "R8" = 561.7673603816401
    "S8" = 0.09930612921153359
    "T8" = 0.0993061292115336
    "E9" = 3
    "G9" = 0.9802360000000001
    "H9" = 2.940708
    "I9" = 0.2336099514851752
    "J9" = 0.2336099514851752
    "K9" = 3
    "M9" = 17.96553866666667
    "N9" = 53.896616
    "O9" = 0.119933805378222
    "P9" = 0.119933805378222
    "Q9" = 17.61046776045867
    "R9" = 158.494209844128
    "S9" = 0.02801773045583889
    "T9" = 0.02801773045583889
    "E10" = 3
    "G10" = 0.8143563333333333
    "H10" = 2.443069
    "I10" = 0.1940774910548533
    "J10" = 0.1940774910548533
    "K10" = 3
    "M10" = 13.441269
    "N10" = 40.323807
    "O10" = 0.08973082133481231
    "P10" = 0.08973082133481232
    "Q10" = 10.945982538187
    "R10" = 98.51384284368301
    "S10" = 0.01741473267495167
    "T10" = 0.01741473267495168
    "E11" = 3
    "G11" = 0.8143563333333333
    "H11" = 2.443069
    "I11" = 0.1940774910548533
    "J11" = 0.1940774910548533
    "K11" = 3
    "M11" = 54.711535
    "N11" = 164.134605
    "O11" = 0.3652416280068742
    "P11" = 0.3652416280068742
    "Q11" = 44.55468503363834
    "R11" = 400.992165302745
    "S11" = 0.07088517879236417
    "T11" = 0.07088517879236418
    "E12" = 3
    "G12" = 0.8143563333333333
    "H12" = 2.443069
    "I12" = 0.1940774910548533
    "J12" = 0.1940774910548533
    "K12" = 3
    "M12" = 63.67711
    "N12" = 191.03133
    "O12" = 0.4250937452800914
    "P12" = 0.4250937452800915
    "Q12" = 51.85585781686333
    "R12" = 466.70272035177
    "S12" = 0.08250112754707102
    "T12" = 0.08250112754707103
    "E13" = 3
    "G13" = 0.8143563333333333
    "H13" = 2.443069
    "I13" = 0.1940774910548533
    "J13" = 0.1940774910548533
    "K13" = 3
    "M13" = 17.96553866666667
    "N13" = 53.896616
    "O13" = 0.119933805378222
    "P13" = 0.119933805378222
    "Q13" = 14.63035019494489
    "R13" = 131.673151754504
    "S13" = 0.02327645204046639
    "T13" = 0.02327645204046639
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}